$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17345.52200362984
$ws.Range("C2").Value = 36844.34838206026
$ws.Range("D2").Value = 71467.80460595629
$ws.Range("E2").Value = 108045.3048360338

$ws.Range("B3").Value = 179041.9261728243
$ws.Range("C3").Value = 331912.1271074681
$ws.Range("D3").Value = 401697.0845889651
$ws.Range("E3").Value = 446585.6229259858

$ws.Range("B4").Value = 19208.01336486229
$ws.Range("C4").Value = 36232.9000462401
$ws.Range("D4").Value = 57445.15061809097
$ws.Range("E4").Value = 74280.8662147133

$ws.Range("B6").Value = 105158.923334647
$ws.Range("C6").Value = 133691.0584462131
$ws.Range("D6").Value = 122732.2392604157
$ws.Range("E6").Value = 100599.9492540442

$ws.Range("B7").Value = 10862.62356895449
$ws.Range("C7").Value = 21728.6595052552
$ws.Range("D7").Value = 23891.34993868103
$ws.Range("E7").Value = 25997.67331348265

$ws.Range("B9").Value = 808858.2267282361
$ws.Range("C9").Value = 1296964.06854493
$ws.Range("D9").Value = 1753991.996364924
$ws.Range("E9").Value = 2160121.135333958

$ws.Range("B12").Value = 784457.6067301839
$ws.Range("C12").Value = 862194.0777817733
$ws.Range("D12").Value = 715143.8465439796
$ws.Range("E12").Value = 529659.8117886288
